$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (columns A:T), reflecting the
# "Natmi following Dr Hou advice" recompute: sending/target clusters now
# iterate over all of ECs/FAPs/sCs (previously ECs was only a target),
# and the derived-specificity / expression numbers were recalculated.
$data = @(
  @("ECs","Sema3c","Plxnd1","ECs",2,0.6666666666666666,0.227933,0.683799,0.005274728560320285,0.005274728560320286,3,1,73.19890333333333,219.59671,0.6596328743217019,0.6596328743217019,16.68444563347667,150.16001070129,0.003479384361510842,0.003479384361510843),
  @("ECs","Sema3c","Plxnd1","FAPs",2,0.6666666666666666,0.227933,0.683799,0.005274728560320285,0.005274728560320286,3,1,9.661646333333335,28.984939,0.08706605224007774,0.08706605224007773,2.202208033695667,19.819872303261,0.000459249792385076,0.000459249792385076),
  @("ECs","Sema3c","Plxnd1","sCs",2,0.6666666666666666,0.227933,0.683799,0.005274728560320285,0.005274728560320286,3,1,28.10860633333333,84.325819,0.2533010734382203,0.2533010734382204,6.406878967375667,57.661910706381,0.001336094406424367,0.001336094406424367),
  @("FAPs","Sema3c","Plxnd1","ECs",3,1,42.071953,126.215859,0.9736112457500787,0.9736112457500787,3,1,73.19890333333333,219.59671,0.6596328743217019,0.6596328743217019,3079.620820691543,27716.58738622389,0.6422259845060573,0.6422259845060573),
  @("FAPs","Sema3c","Plxnd1","FAPs",3,1,42.071953,126.215859,0.9736112457500787,0.9736112457500787,3,1,9.661646333333335,28.984939,0.08706605224007774,0.08706605224007773,406.4843304386224,3658.358973947601,0.08476848758400352,0.0847684875840035),
  @("FAPs","Sema3c","Plxnd1","sCs",3,1,42.071953,126.215859,0.9736112457500787,0.9736112457500787,3,1,28.10860633333333,84.325819,0.2533010734382203,0.2533010734382204,1182.583964551502,10643.25568096352,0.2466167736600179,0.2466167736600179),
  @("sCs","Sema3c","Plxnd1","ECs",3,1,0.912385,2.737155,0.02111402568960099,0.021114025689601,3,1,73.19890333333333,219.59671,0.6596328743217019,0.6596328743217019,66.78558141778333,601.07023276005,0.01392750545413376,0.01392750545413376),
  @("sCs","Sema3c","Plxnd1","FAPs",3,1,0.912385,2.737155,0.02111402568960099,0.021114025689601,3,1,9.661646333333335,28.984939,0.08706605224007774,0.08706605224007773,8.815141189838336,79.33627070854502,0.001838314863689143,0.001838314863689143),
  @("sCs","Sema3c","Plxnd1","sCs",3,1,0.912385,2.737155,0.02111402568960099,0.021114025689601,3,1,28.10860633333333,84.325819,0.2533010734382203,0.2533010734382204,25.64587078943833,230.812837104945,0.005348205371778092,0.005348205371778094)
)

$startRow = 2
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $c).Value = $rowVals[$c - 1]
    }
}
